$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows above row 129 (everything from old row 129
# downward shifts down by two rows, so old row 191/192 become 193/194).
$ws.Rows("129:130").Insert()

# New row 129
$ws.Range("A129").Value = 10
$ws.Range("B129").Value = "Vega Modelo de Temuco"
$ws.Range("C129").Value = "La Araucanía"
$ws.Range("D129").Value = [DateTime]"2021-09-13"
$ws.Range("E129").Value = 9
$ws.Range("F129").Value = 100112037
$ws.Range("G129").Value = "Cebollín"
$ws.Range("H129").Value = "Sin especificar"
$ws.Range("I129").Value = "Primera"
$ws.Range("J129").Value = 40
$ws.Range("K129").Value = 8000
$ws.Range("L129").Value = 8000
$ws.Range("M129").Value = 8000
$ws.Range("N129").Value = "$/docena de paquetes"
$ws.Range("O129").Value = "Provincia de Cautín"
$ws.Range("P129").Value = 667
$ws.Range("Q129").Value = 12
$ws.Range("R129").Value = "Hortaliza"

# New row 130
$ws.Range("A130").Value = 10
$ws.Range("B130").Value = "Vega Modelo de Temuco"
$ws.Range("C130").Value = "La Araucanía"
$ws.Range("D130").Value = [DateTime]"2021-09-13"
$ws.Range("E130").Value = 9
$ws.Range("F130").Value = 100112037
$ws.Range("G130").Value = "Cebollín"
$ws.Range("H130").Value = "Sin especificar"
$ws.Range("I130").Value = "Primera"
$ws.Range("J130").Value = 50
$ws.Range("K130").Value = 5000
$ws.Range("L130").Value = 5000
$ws.Range("M130").Value = 5000
$ws.Range("N130").Value = "$/docena de paquetes"
$ws.Range("O130").Value = "Región de O'Higgins"
$ws.Range("P130").Value = 417
$ws.Range("Q130").Value = 12
$ws.Range("R130").Value = "Hortaliza"
